$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.356646523565985
$ws.Range("D2").Value = 0.1336308149434444
$ws.Range("E2").Value = 0.144697372703579
$ws.Range("F2").Value = 1.883631521800524
$ws.Range("G2").Value = 1.420285920779151
$ws.Range("H2").Value = 1.235672159033498
$ws.Range("I2").Value = 0.9187264827654467
$ws.Range("J2").Value = 0.187217379661881
$ws.Range("L2").Value = 0.7559114855183964
$ws.Range("N2").Value = 1.450459894863499
$ws.Range("B3").Value = 1.281849286799741
$ws.Range("D3").Value = 0.1285500666958228
$ws.Range("E3").Value = 0.1386925272860324
$ws.Range("F3").Value = 1.865972900885978
$ws.Range("G3").Value = 1.383708604102424
$ws.Range("H3").Value = 1.225468251151653
$ws.Range("I3").Value = 0.9328546247175584
$ws.Range("J3").Value = 0.1790023011326483
$ws.Range("L3").Value = 0.7018566950308696
$ws.Range("N3").Value = 1.455169515965252
$ws.Range("B4").Value = 1.236305839650015
$ws.Range("D4").Value = 0.1253930487628878
$ws.Range("E4").Value = 0.1349733344247781
$ws.Range("F4").Value = 1.856423180995932
$ws.Range("G4").Value = 1.362359394279622
$ws.Range("H4").Value = 1.219993312181543
$ws.Range("I4").Value = 0.9420785232702276
$ws.Range("J4").Value = 0.1739229841589847
$ws.Range("L4").Value = 0.6689017199970522
$ws.Range("N4").Value = 1.45863204963095
$ws.Range("B5").Value = 1.217843633096976
$ws.Range("D5").Value = 0.1240970260618823
$ws.Range("E5").Value = 0.133449614835758
$ws.Range("F5").Value = 1.852855594420348
$ws.Range("G5").Value = 1.353936824234495
$ws.Range("H5").Value = 1.217960396098505
$ws.Range("I5").Value = 0.9459752036682616
$ws.Range("J5").Value = 0.1718442954342905
$ws.Range("L5").Value = 0.6555314892683271
$ws.Range("N5").Value = 1.460186956898418
$ws.Range("B6").Value = 1.214783892253962
$ws.Range("D6").Value = 0.1238812449632434
$ws.Range("E6").Value = 0.133196111384418
$ws.Range("F6").Value = 1.852282737246384
$ws.Range("G6").Value = 1.352554972068162
$ws.Range("H6").Value = 1.217634784385211
$ws.Range("I6").Value = 0.9466305647522439
$ws.Range("J6").Value = 0.1714985986650603
$ws.Range("L6").Value = 0.6533149514719696
$ws.Range("N6").Value = 1.460453849893497
$ws.Range("B7").Value = 1.23605645736518
$ws.Range("D7").Value = 0.1253756087960909
$ws.Range("E7").Value = 0.1349528179078767
$ws.Range("F7").Value = 1.85637375682947
$ws.Range("G7").Value = 1.362244683416662
$ws.Range("H7").Value = 1.219965093952595
$ws.Range("I7").Value = 0.9421305173502255
$ws.Range("J7").Value = 0.1738949859259833
$ws.Range("L7").Value = 0.6687211645628963
$ws.Range("N7").Value = 1.458652436541328
$ws.Range("B8").Value = 1.330777632046875
$ws.Range("D8").Value = 0.1318866754996506
$ws.Range("E8").Value = 0.1426335680910427
$ws.Range("F8").Value = 1.877273862804913
$ws.Range("G8").Value = 1.40744284438901
$ws.Range("H8").Value = 1.23198947781421
$ws.Range("I8").Value = 0.9234837535339757
$ws.Range("J8").Value = 0.1843921482146129
$ws.Range("L8").Value = 0.7372246982644981
$ws.Range("N8").Value = 1.451965514126186
$ws.Range("B9").Value = 1.519527337756642
$ws.Range("D9").Value = 0.1443631300921311
$ws.Range("E9").Value = 0.1574417876633376
$ws.Range("F9").Value = 1.92856969965564
$ws.Range("G9").Value = 1.504956355389822
$ws.Range("H9").Value = 1.261868497774003
$ws.Range("I9").Value = 0.8912845181170113
$ws.Range("J9").Value = 0.2046971311774968
$ws.Range("L9").Value = 0.8734259198713517
$ws.Range("N9").Value = 1.443366862261371
$ws.Range("B10").Value = 1.660004228970308
$ws.Range("D10").Value = 0.1533596293290316
$ws.Range("E10").Value = 0.1681701240975713
$ws.Range("F10").Value = 1.9726240045689
$ws.Range("G10").Value = 1.582130046315541
$ws.Range("H10").Value = 1.287704968963112
$ws.Range("I10").Value = 0.8703026817499797
$ws.Range("J10").Value = 0.2194454447633092
$ws.Range("L10").Value = 0.9746451909499001
$ws.Range("N10").Value = 1.439783515511266
$ws.Range("B11").Value = 1.724297514326111
$ws.Range("D11").Value = 0.1574172224438826
$ws.Range("E11").Value = 0.1730186817611497
$ws.Range("F11").Value = 1.994066267473073
$ws.Range("G11").Value = 1.618465275304885
$ws.Range("H11").Value = 1.300312039324723
$ws.Range("I11").Value = 0.8613410875453198
$ws.Range("J11").Value = 0.2261182685991372
$ws.Range("L11").Value = 1.02094666603054
$ws.Range("N11").Value = 1.43874349981931
$ws.Range("B12").Value = 1.748699023034135
$ws.Range("D12").Value = 0.1589488274364896
$ws.Range("E12").Value = 0.1748501733481049
$ws.Range("F12").Value = 2.002388823569021
$ws.Range("G12").Value = 1.632403080041996
$ws.Range("H12").Value = 1.305209538741735
$ws.Range("I12").Value = 0.8580316749139563
$ws.Range("J12").Value = 0.2286398725797341
$ws.Range("L12").Value = 1.038516786327278
$ws.Range("N12").Value = 1.438434227485885
$ws.Range("B13").Value = 1.743441289438238
$ws.Range("D13").Value = 0.1586191858427952
$ws.Range("E13").Value = 0.1744559309357641
$ws.Range("F13").Value = 2.000587370635472
$ws.Range("G13").Value = 1.629393359596577
$ws.Range("H13").Value = 1.304149272049841
$ws.Range("I13").Value = 0.8587406689170685
$ws.Range("J13").Value = 0.2280970340784876
$ws.Range("L13").Value = 1.03473111330166
$ws.Range("N13").Value = 1.438497078573036
$ws.Range("B14").Value = 1.72630394375949
$ws.Range("D14").Value = 0.1575433264287369
$ws.Range("E14").Value = 0.1731694507679507
$ws.Range("F14").Value = 1.994746896708889
$ws.Range("G14").Value = 1.619608360377612
$ws.Range("H14").Value = 1.300712480980195
$ws.Range("I14").Value = 0.861067132028225
$ws.Range("J14").Value = 0.2263258278859723
$ws.Range("L14").Value = 1.022391434291876
$ws.Range("N14").Value = 1.438716362746973
$ws.Range("B15").Value = 1.715813962853701
$ws.Range("D15").Value = 0.1568836943525156
$ws.Range("E15").Value = 0.1723808524238564
$ws.Range("F15").Value = 1.991195890983747
$ws.Range("G15").Value = 1.613638057415898
$ws.Range("H15").Value = 1.298623447755233
$ws.Range("I15").Value = 0.8625031252180513
$ws.Range("J15").Value = 0.2252402282601906
$ws.Range("L15").Value = 1.014837808543263
$ws.Range("N15").Value = 1.438861684220882
$ws.Range("B16").Value = 1.655810266994649
$ws.Range("D16").Value = 0.1530937619474884
$ws.Range("E16").Value = 0.1678526211828881
$ws.Range("F16").Value = 1.971251016886953
$ws.Range("G16").Value = 1.57978032664343
$ws.Range("H16").Value = 1.286898308543158
$ws.Range("I16").Value = 0.8709001023933887
$ws.Range("J16").Value = 0.2190086273474918
$ws.Range("L16").Value = 0.9716244295196645
$ws.Range("N16").Value = 1.43986333146448
$ws.Range("B17").Value = 1.619099029427161
$ws.Range("D17").Value = 0.1507598848340876
$ws.Range("E17").Value = 0.1650665652186589
$ws.Range("F17").Value = 1.959375399064427
$ws.Range("G17").Value = 1.559325622953111
$ws.Range("H17").Value = 1.279924509082292
$ws.Range("I17").Value = 0.8762009415352559
$ws.Range("J17").Value = 0.2151764296862382
$ws.Range("L17").Value = 0.9451799950439295
$ws.Range("N17").Value = 1.44062872910385
$ws.Range("B18").Value = 1.598020430216934
$ws.Range("D18").Value = 0.1494142032537411
$ws.Range("E18").Value = 0.1634611033557576
$ws.Range("F18").Value = 1.952676719355296
$ws.Range("G18").Value = 1.547676145895224
$ws.Range("H18").Value = 1.27599369870353
$ws.Range("I18").Value = 0.8793047150593409
$ws.Range("J18").Value = 0.2129688402352485
$ws.Range("L18").Value = 0.9299939873330345
$ws.Range("N18").Value = 1.441124522298821
$ws.Range("B19").Value = 1.59088991820613
$ws.Range("D19").Value = 0.1489580089780276
$ws.Range("E19").Value = 0.1629170062729131
$ws.Range("F19").Value = 1.950431270445364
$ws.Range("G19").Value = 1.543751621105798
$ws.Range("H19").Value = 1.27467657078472
$ws.Range("I19").Value = 0.8803650147205886
$ws.Range("J19").Value = 0.2122208036959847
$ws.Range("L19").Value = 0.9248564139637381
$ws.Range("N19").Value = 1.441301940318866
$ws.Range("B20").Value = 1.623003208776936
$ws.Range("D20").Value = 0.1510086704930131
$ws.Range("E20").Value = 0.1653634555586052
$ws.Range("F20").Value = 1.960625922549468
$ws.Range("G20").Value = 1.561491091581757
$ws.Range("H20").Value = 1.280658562709789
$ws.Range("I20").Value = 0.8756309771048514
$ws.Range("J20").Value = 0.2155847269274886
$ws.Range("L20").Value = 0.9479925534860172
$ws.Range("N20").Value = 1.440541503455762
$ws.Range("B21").Value = 1.731336111770304
$ws.Range("D21").Value = 0.157859464996136
$ws.Range("E21").Value = 0.17354744452728
$ws.Range("F21").Value = 1.996456870368661
$ws.Range("G21").Value = 1.622477594688547
$ws.Range("H21").Value = 1.301718593257192
$ws.Range("I21").Value = 0.8603815073042611
$ws.Range("J21").Value = 0.226846216830765
$ws.Range("L21").Value = 1.026014900331802
$ws.Range("N21").Value = 1.438649661039037
$ws.Range("B22").Value = 1.802458271187049
$ws.Range("D22").Value = 0.1623082443639134
$ws.Range("E22").Value = 0.1788696335542994
$ws.Range("F22").Value = 2.021057351851368
$ws.Range("G22").Value = 1.663376736651458
$ws.Range("H22").Value = 1.31620257351733
$ws.Range("I22").Value = 0.8509057173384171
$ws.Range("J22").Value = 0.2341756371730952
$ws.Range("L22").Value = 1.077221220135613
$ws.Range("N22").Value = 1.437905978829278
$ws.Range("B23").Value = 1.764470049651891
$ws.Range("D23").Value = 0.1599364284056151
$ws.Range("E23").Value = 0.1760314998350552
$ws.Range("F23").Value = 2.007818954745588
$ws.Range("G23").Value = 1.641452255858297
$ws.Range("H23").Value = 1.308406089880975
$ws.Range("I23").Value = 0.8559181356417547
$ws.Range("J23").Value = 0.2302666002979663
$ws.Range("L23").Value = 1.04987187457715
$ws.Range("N23").Value = 1.438257903626081
$ws.Range("B24").Value = 1.621238044556549
$ws.Range("D24").Value = 0.1508962066667436
$ws.Range("E24").Value = 0.1652292430282074
$ws.Range("F24").Value = 1.960060159884023
$ws.Range("G24").Value = 1.560511740086412
$ws.Range("H24").Value = 1.280326452544557
$ws.Range("I24").Value = 0.8758884828382776
$ws.Range("J24").Value = 0.2154001494719608
$ws.Range("L24").Value = 0.9467209420763822
$ws.Range("N24").Value = 1.440580764481879
$ws.Range("B25").Value = 1.468147065852634
$ws.Range("D25").Value = 0.141018338340885
$ws.Range("E25").Value = 0.1534625845385094
$ws.Range("F25").Value = 1.91358060174511
$ws.Range("G25").Value = 1.477613664473381
$ws.Range("H25").Value = 1.253106522638802
$ws.Range("I25").Value = 0.8995264432972618
$ws.Range("J25").Value = 0.1992339747186946
$ws.Range("L25").Value = 0.8363786342017363
$ws.Range("N25").Value = 1.44521162885485
